$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Swap cell styles: the title row (A1:D1) becomes the "italic" style that
#    used to belong to rows 34/44, and rows 34/44 become the "header/centered"
#    style that used to belong to row 1. We use helper cells far outside the
#    used range to stash formats via copy / paste-special (formats only) so
#    that no new style entries get created - we just re-point existing ones.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("Z100").PasteSpecial(-4122)   # Z100 now holds the header style
$ws.Range("A34").Copy()
$ws.Range("Z101").PasteSpecial(-4122)   # Z101 now holds the italic style

$ws.Range("Z101").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)  # row1 gets the italic style

$ws.Range("Z100").Copy()
$ws.Range("A34:D34").PasteSpecial(-4122)  # row34 gets the header style
$ws.Range("A44:D44").PasteSpecial(-4122)  # row44 gets the header style

$ws.Range("Z100:Z101").Clear()

# ---------------------------------------------------------------------------
# 2) Append the new Assets.* module rows (and update B51, which used to
#    contain "Product.Assets.Data" - that string becomes unused, so it
#    disappears from the shared strings automatically). The B-column
#    (module name) values are written in this specific order so that the
#    shared-strings table ends up with the same ordering as the original
#    edit.
# ---------------------------------------------------------------------------
$ws.Range("A52").Value = 2001
$ws.Range("C52").Value = "A"

$ws.Range("A53").Value = 2002
$ws.Range("C53").Value = "A"

$ws.Range("A54").Value = 2003
$ws.Range("C54").Value = "A"

$ws.Range("A55").Value = 2004
$ws.Range("C55").Value = "A"

$ws.Range("A56").Value = 2005
$ws.Range("C56").Value = "A"

$ws.Range("B53").Value = "Assets.Server"
$ws.Range("B56").Value = "Assets.Server.Engine"
$ws.Range("B55").Value = "Assets.Server.DataFillers"
$ws.Range("B52").Value = "Assets.Data"
$ws.Range("B51").Value = "Assets.Core"
$ws.Range("B54").Value = "Assets.App"

# ---------------------------------------------------------------------------
# 3) Update the selected cell (mirrors the author moving the cursor after
#    appending the new rows).
# ---------------------------------------------------------------------------
$ws.Range("A57").Select() | Out-Null
